# Weekly fruit/vegetable price update: insert a new week's record for
# row 83 (Packham's Triumph, Primera) and shift the remaining Pera cycle
# rows down by one, matching the sheet's repeating weekly-cycle layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 83; Excel shifts rows 83:160 down to
# 84:161 automatically (this also grows the used range to A1:T161).
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with this cycle's new record.
$ws.Range("A83").Value = 7
$ws.Range("B83").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C83").Value = "Ñuble"
$ws.Range("D83").Value = 44512
$ws.Range("E83").Value = 16
$ws.Range("F83").Value = "Fruta"
$ws.Range("G83").Value = 100104
$ws.Range("H83").Value = "Frutos de pepita"
$ws.Range("I83").Value = 100104005
$ws.Range("J83").Value = "Pera"
$ws.Range("K83").Value = "Packham's Triumph"
$ws.Range("L83").Value = "Primera"
$ws.Range("M83").Value = 80
$ws.Range("N83").Value = 8500
$ws.Range("O83").Value = 9000
$ws.Range("P83").Value = 8750
$ws.Range("Q83").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R83").Value = "Provincia de Curicó"
$ws.Range("S83").Value = 547
$ws.Range("T83").Value = 16
